$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Valor Mora total: 120000 -> 180000 ---
$ws.Range("E11").Value = 180000

# --- Cant. Periodos: 2 -> 3 ---
$ws.Range("F13").Value = 3

# --- Swap the "Novedad de Retiro" / "Novedad de Ingreso" header order ---
$ws.Range("H15").Value = "Novedad de Retiro"
$ws.Range("I15").Value = "Novedad de Ingreso"

# --- Add a new "2509" period row, pushing the signature block down ---
$ws.Rows("17:17").Copy()
$ws.Rows("18:18").Insert(-4121)

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1143344035"
$ws.Range("D18").Value = "DARILUZ MARTINEZ GALVIS"
$ws.Range("E18").Value = "2509"
$ws.Range("F18").Value = 60000
$ws.Range("G18").Value = 1500000
$ws.Range("H18").Value = ""
$ws.Range("I18").Value = ""
$ws.Range("J18").Value = ""

Write-Output "done"
